$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "email" column header (K1)
$ws.Range("K1").Value = "email"
$ws.Range("K1").NumberFormat = "@"

# New email value (K2)
$ws.Range("K2").Value = "gb@gmail.com"
$ws.Range("K2").NumberFormat = "@"

# Turn K2 into a real mailto hyperlink, keeping the visible text as the address
$ws.Hyperlinks.Add($ws.Range("K2"), "mailto:gb@gmail.com", "", "", "gb@gmail.com")

# Style it as plain blue text (rather than Excel's default underlined
# "Hyperlink" look) and drop the auto-added named style so it doesn't linger
# unused in the workbook
$ws.Range("K2").Font.Name = "Arial"
$ws.Range("K2").Font.Underline = -4142
$ws.Range("K2").Font.Color = 16711680
foreach ($s in $wb.Styles) {
  if ($s.Name -eq "Hyperlink") {
    $s.Delete()
  }
}

# Give the new column a sensible width
$ws.Columns("K").ColumnWidth = 18.47

# Move/restore the selection (fixes the reported scrolling issue)
$null = $ws.Range("K37").Select()

Write-Host "done"
